# Auto-generated Excel COM-interop script to apply scraped schedule update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = "Última actualización: 07:51:40"
$ws.Cells.Item(3, 1).Value = "Total filas: 58"
$ws.Cells.Item(33, 1).Value = "07:51:40"
$ws.Cells.Item(33, 4).Value = 4
$ws.Cells.Item(35, 1).Value = "07:51:40"
$ws.Cells.Item(35, 3).Value = "17_ROMERO"
$ws.Cells.Item(35, 4).Value = 9
$ws.Cells.Item(36, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(37, 1).Value = "07:51:40"
$ws.Cells.Item(37, 4).Value = 10
$ws.Cells.Item(38, 1).Value = "07:51:40"
$ws.Cells.Item(38, 4).Value = 15
$ws.Cells.Item(39, 1).Value = "07:51:40"
$ws.Cells.Item(39, 4).Value = 20
$ws.Cells.Item(41, 1).Value = "07:51:40"
$ws.Cells.Item(41, 4).Value = 22
$ws.Cells.Item(44, 1).Value = "07:51:40"
$ws.Cells.Item(44, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(44, 4).Value = 38
$ws.Cells.Item(45, 1).Value = "07:51:40"
$ws.Cells.Item(45, 3).Value = "15_ABASTO"
$ws.Cells.Item(45, 4).Value = 38
$ws.Cells.Item(47, 1).Value = "07:51:40"
$ws.Cells.Item(47, 4).Value = 50
$ws.Cells.Item(48, 1).Value = "07:51:40"
$ws.Cells.Item(48, 4).Value = 52
$ws.Cells.Item(49, 1).Value = "07:51:40"
$ws.Cells.Item(49, 2).Value = "08:45"
$ws.Cells.Item(49, 4).Value = 54
$ws.Cells.Item(50, 1).Value = "07:26:49"
$ws.Cells.Item(50, 2).Value = "08:51"
$ws.Cells.Item(50, 4).Value = 85
$ws.Cells.Item(51, 1).Value = "06:58:58"
$ws.Cells.Item(51, 2).Value = "08:52"
$ws.Cells.Item(51, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(51, 4).Value = 114
$ws.Cells.Item(52, 1).Value = "07:51:40"
$ws.Cells.Item(52, 2).Value = "08:53"
$ws.Cells.Item(52, 3).Value = "215B_EL PATO"
$ws.Cells.Item(52, 4).Value = 62
$ws.Cells.Item(53, 2).Value = "08:57"
$ws.Cells.Item(53, 3).Value = "215A_EL PATO"
$ws.Cells.Item(53, 4).Value = 91
$ws.Cells.Item(54, 1).Value = "07:51:40"
$ws.Cells.Item(54, 2).Value = "08:58"
$ws.Cells.Item(54, 3).Value = "215A_EL PATO"
$ws.Cells.Item(54, 4).Value = 67
$ws.Cells.Item(55, 1).Value = "07:51:40"
$ws.Cells.Item(55, 2).Value = "09:06"
$ws.Cells.Item(55, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(55, 4).Value = 75
$ws.Cells.Item(56, 2).Value = "09:16"
$ws.Cells.Item(56, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(56, 4).Value = 110
$ws.Cells.Item(57, 1).Value = "07:51:40"
$ws.Cells.Item(57, 2).Value = "09:17"
$ws.Cells.Item(57, 3).Value = "14_ABASTO"
$ws.Cells.Item(57, 4).Value = 86
$ws.Cells.Item(57, 5).Value = "LP1912"
$ws.Cells.Item(58, 1).Value = "07:51:40"
$ws.Cells.Item(58, 2).Value = "09:18"
$ws.Cells.Item(58, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(58, 4).Value = 87
$ws.Cells.Item(58, 5).Value = "LP1912"
$ws.Cells.Item(59, 1).Value = "07:51:40"
$ws.Cells.Item(59, 2).Value = "09:21"
$ws.Cells.Item(59, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(59, 4).Value = 90
$ws.Cells.Item(59, 5).Value = "LP1912"
$ws.Cells.Item(60, 1).Value = "07:51:40"
$ws.Cells.Item(60, 2).Value = "09:29"
$ws.Cells.Item(60, 3).Value = "10_OLMOS"
$ws.Cells.Item(60, 4).Value = 98
$ws.Cells.Item(60, 5).Value = "LP1912"
$ws.Cells.Item(61, 1).Value = "07:51:40"
$ws.Cells.Item(61, 2).Value = "09:39"
$ws.Cells.Item(61, 3).Value = "15_ABASTO"
$ws.Cells.Item(61, 4).Value = 108
$ws.Cells.Item(61, 5).Value = "LP1912"
$ws.Cells.Item(62, 1).Value = "07:51:40"
$ws.Cells.Item(62, 2).Value = "09:41"
$ws.Cells.Item(62, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(62, 4).Value = 110
$ws.Cells.Item(62, 5).Value = "LP1912"
$ws.Cells.Item(63, 1).Value = "07:51:40"
$ws.Cells.Item(63, 2).Value = "09:43"
$ws.Cells.Item(63, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(63, 4).Value = 112
$ws.Cells.Item(63, 5).Value = "LP1912"

$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = "Última actualización: 07:51:40"
$ws.Cells.Item(3, 1).Value = "Total filas: 11"
$ws.Cells.Item(13, 1).Value = "07:51:40"
$ws.Cells.Item(13, 4).Value = 52
$ws.Cells.Item(14, 1).Value = "07:51:40"
$ws.Cells.Item(14, 4).Value = 62
$ws.Cells.Item(16, 1).Value = "07:51:40"
$ws.Cells.Item(16, 2).Value = "08:58"
$ws.Cells.Item(16, 3).Value = "215A_EL PATO"
$ws.Cells.Item(16, 4).Value = 67
$ws.Cells.Item(16, 5).Value = "LP1912"

$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = "Última actualización: 07:51:40"
$ws.Cells.Item(3, 1).Value = "Total filas: 5"
$ws.Cells.Item(8, 1).Value = "07:51:40"
$ws.Cells.Item(8, 4).Value = 44
$ws.Cells.Item(9, 1).Value = "07:51:40"
$ws.Cells.Item(9, 4).Value = 59
$ws.Cells.Item(10, 1).Value = "07:51:40"
$ws.Cells.Item(10, 2).Value = "09:20"
$ws.Cells.Item(10, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(10, 4).Value = 89
$ws.Cells.Item(10, 5).Value = "L6173"
